$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("A4").Value = 131130604
$ws.Range("B4").Value = 57073
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 100138
$ws.Range("F4").Value = "Tjäder"
$ws.Range("G4").Value = "Tetrao urogallus"
$ws.Range("H4").Value = "Linnaeus, 1758"
$ws.Range("I4").Value = "'3"
$ws.Range("K4").Value = "adult"
$ws.Range("L4").Value = "hona"
$ws.Range("M4").Value = "födosökande"
$ws.Range("P4").Value = "Tågmossen, Tågmossen, Vstm"
$ws.Range("Q4").Value = 550930
$ws.Range("R4").Value = 6620382
$ws.Range("S4").Value = 50
$ws.Range("T4").Value = "Västmanland"
$ws.Range("U4").Value = "Skinnskatteberg"
$ws.Range("V4").Value = "Västmanland"
$ws.Range("W4").Value = "Gunnilbo"
$ws.Range("Y4").Value = "'2026-02-12"
$ws.Range("Z4").Value = "15:28"
$ws.Range("AA4").Value = "'2026-02-12"
$ws.Range("AB4").Value = "15:28"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AT4").Value = "'"
$ws.Range("AW4").Value = "Per Eriksson"
$ws.Range("AX4").Value = "Per Eriksson, Lena Öling"
$ws.Range("AY4").Value = "'"
